$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.902.02"
$ws.Cells.Item(2, 5).Value = "  -0.06%  "
$ws.Cells.Item(3, 4).Value = "1.633.86"
$ws.Cells.Item(3, 5).Value = "  -0.86%  "
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "212.09"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.76%  "
$ws.Cells.Item(7, 5).Value = "  -0.05%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "23.20"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.46%  "
$ws.Cells.Item(9, 5).Value = "  -3.30%  "
$ws.Cells.Item(10, 5).Value = "  -0.38%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0880"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "1.866.07"
$ws.Cells.Item(12, 5).Value = "  -0.82%  "
$ws.Cells.Item(13, 4).Value = "1.650.45"
$ws.Cells.Item(13, 5).Value = "  +0.17%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.568"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.69%  "
$ws.Cells.Item(16, 5).Value = "  -0.80%  "
$ws.Cells.Item(17, 4).Value = "27.901.54"
$ws.Cells.Item(17, 5).Value = "  -0.07%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "229.44"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -1.25%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0720"
$ws.Cells.Item(19, 5).Value = "  -0.47%  "
$ws.Cells.Item(20, 5).Value = "  -2.25%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.00"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.08%  "
$ws.Cells.Item(22, 5).Value = "  -0.79%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.35"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -3.63%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.07"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -3.99%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "153.04"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.47%  "
$ws.Cells.Item(26, 5).Value = "  +0.64%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "15.61"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.70%  "
$ws.Cells.Item(28, 5).Value = "  -0.67%  "
$ws.Cells.Item(29, 5).Value = "  -0.01%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.19"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.89%  "
$ws.Cells.Item(31, 5).Value = "  -0.68%  "
$ws.Cells.Item(32, 5).Value = "  +0.63%  "
$ws.Cells.Item(33, 4).Value = "1.400.29"
$ws.Cells.Item(33, 5).Value = "  -3.55%  "
$ws.Cells.Item(34, 5).Value = "  -1.64%  "
$ws.Cells.Item(35, 5).Value = "  +0.82%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.00"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +8.53%  "
$ws.Cells.Item(37, 5).Value = "  +1.57%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0169"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.43%  "
$ws.Cells.Item(39, 5).Value = "  -0.20%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.871"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.27%  "
$ws.Cells.Item(42, 5).Value = "  -0.06%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "66.85"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -3.59%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.50"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +2.24%  "
$ws.Cells.Item(45, 5).Value = "  +0.81%  "
$ws.Cells.Item(46, 5).Value = "  -1.73%  "
$ws.Cells.Item(47, 4).Value = "1.775.43"
$ws.Cells.Item(47, 5).Value = "  -0.83%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "87.64"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.30%  "
$ws.Cells.Item(49, 5).Value = "  -0.55%  "
$ws.Cells.Item(50, 5).Value = "  -0.22%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "7.55"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.82%  "
